# Auto-generated edit script applying the Leve profit-sheet value updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 511.25
$ws.Range("I6").Value = 460
$ws.Range("J6").Value = 562.5
$ws.Range("K6").Value = 1380
$ws.Range("L6").Value = 1687.5
$ws.Range("M6").Value = -1268
$ws.Range("N6").Value = -1911.5

$ws.Range("H32").Value = 7695115.5
$ws.Range("J32").Value = 10003290
$ws.Range("L32").Value = 10003290
$ws.Range("N32").Value = -10003942

$ws.Range("H43").Value = 4677.3
$ws.Range("J43").Value = 4599.6
$ws.Range("L43").Value = 4599.6
$ws.Range("N43").Value = -4737.6

$ws.Range("H55").Value = 513.7857
$ws.Range("J55").Value = 405.25
$ws.Range("L55").Value = 405.25
$ws.Range("N55").Value = -833.25

$ws.Range("H92").Value = 1744.4706
$ws.Range("I92").Value = 1853.9231
$ws.Range("K92").Value = 1853.9231
$ws.Range("M92").Value = -605.9231

$ws.Range("H100").Value = 6443.44
$ws.Range("I100").Value = 5473.9165
$ws.Range("K100").Value = 5473.9165
$ws.Range("M100").Value = -4932.9165

$ws.Range("H112").Value = 1790.3462
$ws.Range("I112").Value = 1182.5
$ws.Range("J112").Value = 1841
$ws.Range("K112").Value = 3547.5
$ws.Range("L112").Value = 5523
$ws.Range("M112").Value = -2439.5
$ws.Range("N112").Value = -7739

$ws.Range("H132").Value = 15139.225
$ws.Range("I132").Value = 2238.4849
$ws.Range("K132").Value = 6715.4547
$ws.Range("M132").Value = -4185.4547

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4990.074
$ws.Range("I32").Value = 5307.9185
$ws.Range("J32").Value = 1875.2
$ws.Range("K32").Value = 5307.9185
$ws.Range("L32").Value = 1875.2
$ws.Range("M32").Value = -5020.9185
$ws.Range("N32").Value = -2449.2

$ws.Range("H61").Value = 24332.666
$ws.Range("I61").Value = 24332.666
$ws.Range("K61").Value = 24332.666
$ws.Range("M61").Value = -24120.666

$ws.Range("H102").Value = 3378
$ws.Range("I102").Value = 3635.9333
$ws.Range("J102").Value = 1443.5
$ws.Range("K102").Value = 3635.9333
$ws.Range("L102").Value = 1443.5
$ws.Range("M102").Value = -2013.9333
$ws.Range("N102").Value = -4687.5

$ws.Range("H122").Value = 3281.5557
$ws.Range("I122").Value = 3047.6086
$ws.Range("J122").Value = 4626.75
$ws.Range("K122").Value = 9142.825800000001
$ws.Range("L122").Value = 13880.25
$ws.Range("M122").Value = -6692.825800000001
$ws.Range("N122").Value = -18780.25

$ws.Range("H136").Value = 24332.666
$ws.Range("I136").Value = 24332.666
$ws.Range("K136").Value = 72997.99800000001
$ws.Range("M136").Value = -70447.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1209.8
$ws.Range("I5").Value = 683.3333
$ws.Range("K5").Value = 683.3333
$ws.Range("M5").Value = -570.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 200
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -424

$ws.Range("H31").Value = 1431.3721
$ws.Range("I31").Value = 1346.4048
$ws.Range("J31").Value = 5000
$ws.Range("K31").Value = 1346.4048
$ws.Range("L31").Value = 5000
$ws.Range("M31").Value = -1051.4048
$ws.Range("N31").Value = -5590

$ws.Range("H34").Value = 1431.3721
$ws.Range("I34").Value = 1346.4048
$ws.Range("J34").Value = 5000
$ws.Range("K34").Value = 1346.4048
$ws.Range("L34").Value = 5000
$ws.Range("M34").Value = -1144.4048
$ws.Range("N34").Value = -5404

$ws.Range("H132").Value = 2356.2856
$ws.Range("I132").Value = 1682.5
$ws.Range("J132").Value = 6399
$ws.Range("K132").Value = 5047.5
$ws.Range("L132").Value = 19197
$ws.Range("M132").Value = -2517.5
$ws.Range("N132").Value = -24257

$ws.Range("H134").Value = 2499
$ws.Range("I134").Value = 2499
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7497
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -4962
$ws.Range("N134").ClearContents()

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 8325.23
$ws.Range("I7").Value = 8648.166999999999
$ws.Range("J7").Value = 4450
$ws.Range("K7").Value = 25944.501
$ws.Range("L7").Value = 13350
$ws.Range("M7").Value = -25832.501
$ws.Range("N7").Value = -13574

$ws.Range("H12").Value = 8247.200000000001
$ws.Range("I12").Value = 1380
$ws.Range("J12").Value = 11680.8
$ws.Range("K12").Value = 4140
$ws.Range("L12").Value = 35042.39999999999
$ws.Range("M12").Value = -3967
$ws.Range("N12").Value = -35388.39999999999

$ws.Range("H86").Value = 1231.9166
$ws.Range("I86").Value = 1569
$ws.Range("K86").Value = 4707
$ws.Range("M86").Value = -3521

$ws.Range("H89").Value = 1231.9166
$ws.Range("I89").Value = 1569
$ws.Range("K89").Value = 14121
$ws.Range("M89").Value = -8193

$ws.Range("H99").Value = 71319.8
$ws.Range("I99").Value = 2383.3333
$ws.Range("J99").Value = 174724.5
$ws.Range("K99").Value = 7149.999899999999
$ws.Range("L99").Value = 524173.5
$ws.Range("M99").Value = -4903.999899999999
$ws.Range("N99").Value = -528665.5

$ws.Range("H107").Value = 2632.4546
$ws.Range("J107").Value = 2581.8696
$ws.Range("L107").Value = 7745.6088
$ws.Range("N107").Value = -11585.6088

$ws.Range("H131").Value = 3520.8462
$ws.Range("J131").Value = 6403.4
$ws.Range("L131").Value = 19210.2
$ws.Range("N131").Value = -29290.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10004
$ws.Range("I5").Value = 10004
$ws.Range("K5").Value = 10004
$ws.Range("M5").Value = -9892

$ws.Range("H49").Value = 36495
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()

$ws.Range("H80").Value = 6287.533
$ws.Range("I80").Value = 3909.25
$ws.Range("K80").Value = 3909.25
$ws.Range("M80").Value = -2911.25

$ws.Range("H83").Value = 6287.533
$ws.Range("I83").Value = 3909.25
$ws.Range("K83").Value = 19546.25
$ws.Range("M83").Value = -14554.25

$ws.Range("H122").Value = 3046.8333
$ws.Range("I122").Value = 2256.7144
$ws.Range("J122").Value = 5812.25
$ws.Range("K122").Value = 6770.1432
$ws.Range("L122").Value = 17436.75
$ws.Range("M122").Value = -4320.1432
$ws.Range("N122").Value = -22336.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3561.4
$ws.Range("I7").Value = 2727
$ws.Range("J7").Value = 6899
$ws.Range("K7").Value = 2727
$ws.Range("L7").Value = 6899
$ws.Range("M7").Value = -2615
$ws.Range("N7").Value = -7123

$ws.Range("H42").Value = 27250
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -31126

$ws.Range("H46").Value = 3139.7334
$ws.Range("J46").Value = 3599.6667
$ws.Range("L46").Value = 3599.6667
$ws.Range("N46").Value = -3975.6667

$ws.Range("H49").Value = 27250
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30294

$ws.Range("H68").Value = 10002
$ws.Range("I68").Value = 10002
$ws.Range("K68").Value = 10002
$ws.Range("M68").Value = -9253

$ws.Range("H71").Value = 10002
$ws.Range("I71").Value = 10002
$ws.Range("K71").Value = 50010
$ws.Range("M71").Value = -46266

$ws.Range("H126").Value = 3561.4
$ws.Range("I126").Value = 2727
$ws.Range("J126").Value = 6899
$ws.Range("K126").Value = 8181
$ws.Range("L126").Value = 20697
$ws.Range("M126").Value = -5711
$ws.Range("N126").Value = -25637

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11132278
$ws.Range("J2").Value = 20071.428
$ws.Range("L2").Value = 20071.428
$ws.Range("N2").Value = -20295.428

$ws.Range("H81").Value = 7410204.5
$ws.Range("I81").Value = 1505.6
$ws.Range("J81").Value = 22227602
$ws.Range("K81").Value = 3011.2
$ws.Range("L81").Value = 44455204
$ws.Range("M81").Value = -1950.2
$ws.Range("N81").Value = -44457326

$ws.Range("H84").Value = 7410204.5
$ws.Range("I84").Value = 1505.6
$ws.Range("J84").Value = 22227602
$ws.Range("K84").Value = 15056
$ws.Range("L84").Value = 222276020
$ws.Range("M84").Value = -9752
$ws.Range("N84").Value = -222286628

$ws.Range("H107").Value = 340.8889
$ws.Range("I107").Value = 340.8889
$ws.Range("K107").Value = 1022.6667
$ws.Range("M107").Value = 897.3333

$ws.Range("H125").Value = 42058.824
$ws.Range("J125").Value = 40812.5
$ws.Range("L125").Value = 40812.5
$ws.Range("N125").Value = -50652.5

$ws.Range("H132").Value = 3800.875
$ws.Range("I132").Value = 3481.7742
$ws.Range("K132").Value = 10445.3226
$ws.Range("M132").Value = -7915.3226

$ws.Range("H136").Value = 1172.5454
$ws.Range("I136").Value = 1139.8
$ws.Range("K136").Value = 3419.4
$ws.Range("M136").Value = -869.3999999999996
